# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data using K instead of Strike#.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 3
    9  = 1
    10 = 2
    11 = 2
    12 = 0
    13 = 0
    15 = 0
    16 = 1
    17 = 2
    18 = 0
    19 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 2
    25 = 2
    26 = 1
    27 = 0
    29 = 0
    31 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
